$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D sometimes hold numeric-looking text (e.g. "574.41").
# Excel auto-converts such strings to real numbers on assignment, which
# would lose the original text formatting (trailing zeros, thousand-dot
# grouping, etc). Force those specific cells to Text format first, then
# restore the default "Normal" style so no stray formatting is left behind.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "61.058.11"
$ws.Range("E2").Value = "  +0.71%  "
$ws.Range("D3").Value = "2.655.61"
$ws.Range("E3").Value = "  +0.97%  "
$ws.Range("E4").Value = "  -0.05%  "
Set-TextValue $ws.Range("D5") "574.41"
$ws.Range("E5").Value = "  -0.36%  "
Set-TextValue $ws.Range("D6") "144.71"
$ws.Range("E6").Value = "  +1.16%  "
Set-TextValue $ws.Range("D7") "0.997"
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  -0.14%  "
Set-TextValue $ws.Range("D9") "6.60"
$ws.Range("E9").Value = "  +1.55%  "
$ws.Range("E10").Value = "  +0.89%  "
$ws.Range("E11").Value = "  +3.92%  "
$ws.Range("E12").Value = "  +1.11%  "
$ws.Range("D13").Value = "3.127.33"
$ws.Range("E13").Value = "  +1.21%  "
Set-TextValue $ws.Range("D14") "26.15"
$ws.Range("E14").Value = "  +12.17%  "
$ws.Range("D15").Value = "61.053.28"
$ws.Range("E15").Value = "  +0.31%  "
$ws.Range("E16").Value = "  +1.00%  "
$ws.Range("D17").Value = "2.664.99"
$ws.Range("E17").Value = "  +1.46%  "
$ws.Range("E18").Value = "  +3.47%  "
$ws.Range("E19").Value = "  +1.60%  "
Set-TextValue $ws.Range("D20") "350.90"
$ws.Range("E20").Value = "  +0.62%  "
$ws.Range("E21").Value = "  +0.87%  "
$ws.Range("E22").Value = "  +0.28%  "
$ws.Range("E23").Value = "  +2.13%  "
Set-TextValue $ws.Range("D24") "64.22"
$ws.Range("E24").Value = "  +1.50%  "
$ws.Range("E25").Value = "  +0.69%  "
$ws.Range("E26").Value = "  -0.12%  "
$ws.Range("E27").Value = "  +5.34%  "
Set-TextValue $ws.Range("D28") "1.98"
$ws.Range("E28").Value = "  +7.51%  "
$ws.Range("D29").Value = "0.0₃0814"
$ws.Range("E29").Value = "  +2.57%  "
Set-TextValue $ws.Range("D30") "6.89"
$ws.Range("E31").Value = "  -0.07%  "
Set-TextValue $ws.Range("D32") "166.50"
$ws.Range("E32").Value = "  +2.60%  "
Set-TextValue $ws.Range("D33") "19.95"
$ws.Range("E33").Value = "  +2.03%  "
$ws.Range("E34").Value = "  +6.94%  "
$ws.Range("E35").Value = "  +8.62%  "
$ws.Range("E36").Value = "  +7.44%  "
$ws.Range("E37").Value = "  +4.58%  "
Set-TextValue $ws.Range("D38") "338.54"
$ws.Range("E38").Value = "  +12.57%  "
$ws.Range("E39").Value = "  +4.08%  "
Set-TextValue $ws.Range("D40") "0.901"
$ws.Range("E40").Value = "  +6.42%  "
Set-TextValue $ws.Range("D41") "38.53"
$ws.Range("E41").Value = "  +1.66%  "
Set-TextValue $ws.Range("D42") "5.26"
$ws.Range("E42").Value = "  +5.63%  "
Set-TextValue $ws.Range("D43") "20.45"
$ws.Range("E43").Value = "  +2.62%  "
Set-TextValue $ws.Range("D44") "134.07"
$ws.Range("E44").Value = "  -0.32%  "
$ws.Range("E45").Value = "  +2.76%  "
Set-TextValue $ws.Range("D46") "0.0998"
$ws.Range("E46").Value = "  +1.30%  "
$ws.Range("E47").Value = "  +2.76%  "
Set-TextValue $ws.Range("D48") "0.616"
$ws.Range("E48").Value = "  +1.56%  "
Set-TextValue $ws.Range("D49") "20.58"
$ws.Range("E50").Value = "  +0.03%  "
$ws.Range("D51").Value = "2.103.74"
$ws.Range("E51").Value = "  +3.88%  "
